# Add a new "2023" column (column L) to the 9.5.1 data-excel worksheet,
# mirroring the existing 2022 column (K) in formatting, and bump the
# height of row 5 (the data row) to fit the now-taller header wrapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (the thin bottom-border spacer row above the header) ---
# K3 is an empty, bordered cell; clone its formatting into L3.
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# --- Row 4 (the year header row) ---
# K4 holds the "2022" header; clone its formatting into L4, then set 2023.
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 2023

# --- Row 5 (the data row with the percentage values) ---
# K5 holds the 2022 value; clone its formatting into L5, then set the
# new 2023 figure.
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 0.11972285283622097

# Row 5 grows a bit taller (wrapped header text needs more room).
$ws.Rows.Item(5).RowHeight = 40.5

# Clear clipboard marching-ants state left over from the copy operations.
$excel.CutCopyMode = 0
